$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range('D2') '25.778.83'
$ws.Range('E2').Value = '  +5.58%  '

# Row 3
Set-TextValue $ws.Range('D3') '1.704.41'
$ws.Range('E3').Value = '  +3.44%  '

# Row 4
$ws.Range('E4').Value = '  -0.06%  '

# Row 5
Set-TextValue $ws.Range('D5') '331.04'
$ws.Range('E5').Value = '  +6.51%  '

# Row 6
Set-TextValue $ws.Range('D6') '0.9994'
$ws.Range('E6').Value = '  +0.09%  '

# Row 7
Set-TextValue $ws.Range('D7') '0.3688'
$ws.Range('E7').Value = '  +1.15%  '

# Row 8
Set-TextValue $ws.Range('D8') '48.17'
$ws.Range('E8').Value = '  +3.03%  '

# Row 9
Set-TextValue $ws.Range('D9') '0.3311'
$ws.Range('E9').Value = '  +2.03%  '

# Row 10
Set-TextValue $ws.Range('D10') '1.168'
$ws.Range('E10').Value = '  +4.08%  '

# Row 11
$ws.Range('E11').Value = '  +4.62%  '

# Row 12
$ws.Range('E12').Value = '  +0.06%  '

# Row 13
Set-TextValue $ws.Range('D13') '6.189'
$ws.Range('E13').Value = '  +3.91%  '

# Row 14
Set-TextValue $ws.Range('D14') '20.01'
$ws.Range('E14').Value = '  +3.27%  '

# Row 15
Set-TextValue $ws.Range('D15') '6.859'
$ws.Range('E15').Value = '  +4.12%  '

# Row 16
Set-TextValue $ws.Range('D16') '1.702.67'
$ws.Range('E16').Value = '  +3.51%  '

# Row 17
Set-TextValue $ws.Range('D17') '0.00001066'
$ws.Range('E17').Value = '  +2.81%  '

# Row 18
Set-TextValue $ws.Range('D18') '0.06623'
$ws.Range('E18').Value = '  +0.90%  '

# Row 19
Set-TextValue $ws.Range('D19') '81.09'
$ws.Range('E19').Value = '  +3.24%  '

# Row 20
Set-TextValue $ws.Range('D20') '0.9993'
$ws.Range('E20').Value = '  -0.05%  '

# Row 21
$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue $ws.Range('D21') '16.16'
$ws.Range('E21').Value = '  +3.56%  '

# Row 22
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue $ws.Range('D22') '6.045'
$ws.Range('E22').Value = '  +1.98%  '

# Row 23
Set-TextValue $ws.Range('D23') '12.95'
$ws.Range('E23').Value = '  +3.32%  '

# Row 24
Set-TextValue $ws.Range('D24') '25.748.72'
$ws.Range('E24').Value = '  +5.40%  '

# Row 25
Set-TextValue $ws.Range('D25') '2.465'
$ws.Range('E25').Value = '  +0.35%  '

# Row 26
Set-TextValue $ws.Range('D26') '2.476'
$ws.Range('E26').Value = '  +6.75%  '

# Row 27
Set-TextValue $ws.Range('D27') '149.42'
$ws.Range('E27').Value = '  +2.09%  '

# Row 28
Set-TextValue $ws.Range('D28') '19.15'
$ws.Range('E28').Value = '  +3.10%  '

# Row 29
Set-TextValue $ws.Range('D29') '1.298'
$ws.Range('E29').Value = '  +9.30%  '

# Row 30
Set-TextValue $ws.Range('D30') '1.891.92'
$ws.Range('E30').Value = '  +3.55%  '

# Row 31
Set-TextValue $ws.Range('D31') '127.72'
$ws.Range('E31').Value = '  +2.99%  '

# Row 32
Set-TextValue $ws.Range('D32') '4.099'
$ws.Range('E32').Value = '  +0.72%  '

# Row 33
Set-TextValue $ws.Range('D33') '5.941'
$ws.Range('E33').Value = '  +4.82%  '

# Row 34
Set-TextValue $ws.Range('D34') '0.08501'
$ws.Range('E34').Value = '  +1.05%  '

# Row 35
Set-TextValue $ws.Range('D35') '1.695'
$ws.Range('E35').Value = '  +3.09%  '

# Row 36
Set-TextValue $ws.Range('D36') '12.82'
$ws.Range('E36').Value = '  +6.15%  '

# Row 37
Set-TextValue $ws.Range('D37') '5.316'
$ws.Range('E37').Value = '  +2.57%  '

# Row 38
Set-TextValue $ws.Range('D38') '1.274'
$ws.Range('E38').Value = '  +0.71%  '

# Row 39
Set-TextValue $ws.Range('D39') '0.06190'
$ws.Range('E39').Value = '  +3.05%  '

# Row 40
Set-TextValue $ws.Range('D40') '8.534'
$ws.Range('E40').Value = '  +5.37%  '

# Row 41
Set-TextValue $ws.Range('D41') '0.2117'
$ws.Range('E41').Value = '  +2.97%  '

# Row 42
Set-TextValue $ws.Range('D42') '0.02245'
$ws.Range('E42').Value = '  +1.03%  '

# Row 43
Set-TextValue $ws.Range('D43') '14.80'
$ws.Range('E43').Value = '  +18.51%  '

# Row 44
Set-TextValue $ws.Range('D44') '0.6099'
$ws.Range('E44').Value = '  +3.65%  '

# Row 45
Set-TextValue $ws.Range('D45') '0.9997'
$ws.Range('E45').Value = '  -0.07%  '

# Row 46
Set-TextValue $ws.Range('D46') '3.841'
$ws.Range('E46').Value = '  +2.33%  '

# Row 47
Set-TextValue $ws.Range('D47') '0.5827'
$ws.Range('E47').Value = '  +3.99%  '

# Row 48
Set-TextValue $ws.Range('D48') '126.62'
$ws.Range('E48').Value = '  +3.58%  '

# Row 49
$ws.Range('E49').Value = '  +3.01%  '

# Row 50
Set-TextValue $ws.Range('D50') '0.07219'
$ws.Range('E50').Value = '  +4.54%  '

# Row 51
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range('D51') '76.42'
$ws.Range('E51').Value = '  +2.91%  '
